$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.532.41'
$ws.Range("E2").Value = '  -1.10%  '
$ws.Range("D3").Value = '1.594.57'
$ws.Range("E3").Value = '  -1.94%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.91'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.504'
$ws.Range("E6").Value = '  -3.48%  '
$ws.Range("E7").Value = '  +0.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.26'
$ws.Range("E8").Value = '  -4.31%  '
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("E10").Value = '  -3.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0866'
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("D12").Value = '1.820.96'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '1.605.07'
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("E14").Value = '  -3.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.538'
$ws.Range("E15").Value = '  -3.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.29'
$ws.Range("E16").Value = '  -3.08%  '
$ws.Range("D17").Value = '27.533.53'
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '218.44'
$ws.Range("E18").Value = '  -4.77%  '
$ws.Range("E19").Value = '  -3.29%  '
$ws.Range("E20").Value = '  -3.60%  '
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.21'
$ws.Range("E22").Value = '  -2.43%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.66'
$ws.Range("E23").Value = '  -4.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.98'
$ws.Range("E24").Value = '  -2.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.62'
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.76'
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("E27").Value = '  +0.42%  '
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("E29").Value = '  -4.04%  '
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("E31").Value = '  -2.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.26'
$ws.Range("E32").Value = '  -4.43%  '
$ws.Range("D33").Value = '1.360.66'
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.96'
$ws.Range("E34").Value = '  -4.17%  '
$ws.Range("E35").Value = '  -2.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.966'
$ws.Range("E36").Value = '  -4.24%  '
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0166'
$ws.Range("E38").Value = '  -2.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.540'
$ws.Range("E39").Value = '  -2.40%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.814'
$ws.Range("E41").Value = '  +0.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.971'
$ws.Range("E42").Value = '  -3.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.37'
$ws.Range("E43").Value = '  -1.06%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.78'
$ws.Range("E44").Value = '  -2.99%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.00'
$ws.Range("E45").Value = '  -2.42%  '
$ws.Range("D46").Value = '1.731.72'
$ws.Range("E46").Value = '  -1.91%  '
$ws.Range("E47").Value = '  -2.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '87.72'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0971'
$ws.Range("E50").Value = '  -3.77%  '
$ws.Range("E51").Value = '  -1.07%  '
